$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5310
$ws.Range("I43").Value = 5931.391
$ws.Range("J43").Value = 4469.294
$ws.Range("K43").Value = 5931.391
$ws.Range("L43").Value = 4469.294
$ws.Range("M43").Value = -5862.391
$ws.Range("N43").Value = -4607.294
$ws.Range("H51").Value = 8378.625
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8378.625
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8378.625
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9346.625
$ws.Range("H53").Value = 336.625
$ws.Range("I53").Value = 282.5
$ws.Range("J53").Value = 499
$ws.Range("K53").Value = 282.5
$ws.Range("L53").Value = 499
$ws.Range("M53").Value = 354.5
$ws.Range("N53").Value = -1773
$ws.Range("H70").Value = 2310.6667
$ws.Range("J70").Value = 2274.3333
$ws.Range("L70").Value = 6822.999899999999
$ws.Range("N70").Value = -7362.999899999999
$ws.Range("H73").Value = 2310.6667
$ws.Range("J73").Value = 2274.3333
$ws.Range("L73").Value = 6822.999899999999
$ws.Range("N73").Value = -8694.999899999999
$ws.Range("H96").Value = 621.0526
$ws.Range("I96").Value = 391.6154
$ws.Range("J96").Value = 1118.1666
$ws.Range("K96").Value = 1174.8462
$ws.Range("L96").Value = 3354.4998
$ws.Range("M96").Value = 198.1538
$ws.Range("N96").Value = -6100.4998
$ws.Range("H135").Value = 2781.5557
$ws.Range("I135").Value = 841.8333
$ws.Range("J135").Value = 6661
$ws.Range("K135").Value = 7576.4997
$ws.Range("L135").Value = 59949
$ws.Range("M135").Value = -5041.4997
$ws.Range("N135").Value = -65019

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6945.2144
$ws.Range("I61").Value = 6244.5835
$ws.Range("K61").Value = 6244.5835
$ws.Range("M61").Value = -6032.5835
$ws.Range("H102").Value = 492.5
$ws.Range("I102").Value = 492.5
$ws.Range("K102").Value = 492.5
$ws.Range("M102").Value = 1129.5
$ws.Range("H132").Value = 2250.0356
$ws.Range("I132").Value = 1492.091
$ws.Range("K132").Value = 4476.272999999999
$ws.Range("M132").Value = -1946.272999999999
$ws.Range("H136").Value = 6945.2144
$ws.Range("I136").Value = 6244.5835
$ws.Range("K136").Value = 18733.7505
$ws.Range("M136").Value = -16183.7505

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -4711
$ws.Range("N29").ClearContents()
$ws.Range("H134").Value = 2777.3635
$ws.Range("I134").Value = 1793.3572
$ws.Range("J134").Value = 8287.799999999999
$ws.Range("K134").Value = 5380.071599999999
$ws.Range("L134").Value = 24863.4
$ws.Range("M134").Value = -2845.071599999999
$ws.Range("N134").Value = -29933.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 47899.855
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 47899.855
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 47899.855
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -49371.855
$ws.Range("H59").Value = 69099.45
$ws.Range("J59").Value = 75566.11
$ws.Range("L59").Value = 75566.11
$ws.Range("N59").Value = -77856.11
$ws.Range("H60").Value = 21309.8
$ws.Range("I60").Value = 9499.666999999999
$ws.Range("J60").Value = 26371.285
$ws.Range("K60").Value = 9499.666999999999
$ws.Range("L60").Value = 26371.285
$ws.Range("M60").Value = -8988.666999999999
$ws.Range("N60").Value = -27393.285
$ws.Range("H61").Value = 47899.855
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 47899.855
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 47899.855
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -48595.855
$ws.Range("H86").Value = 7999
$ws.Range("I86").Value = 7998
$ws.Range("K86").Value = 7998
$ws.Range("M86").Value = -6875
$ws.Range("H89").Value = 7999
$ws.Range("I89").Value = 7998
$ws.Range("K89").Value = 39990
$ws.Range("M89").Value = -34374
$ws.Range("H98").Value = 88000
$ws.Range("J98").Value = 88000
$ws.Range("L98").Value = 88000
$ws.Range("N98").Value = -92492
$ws.Range("H132").Value = 4169.8667
$ws.Range("I132").Value = 2639.55
$ws.Range("K132").Value = 7918.650000000001
$ws.Range("M132").Value = -5388.650000000001
$ws.Range("H134").Value = 4556.914
$ws.Range("I134").Value = 2864.16
$ws.Range("K134").Value = 8592.48
$ws.Range("M134").Value = -6057.48
$ws.Range("H138").Value = 51166.5
$ws.Range("I138").Value = 33333
$ws.Range("K138").Value = 33333
$ws.Range("M138").Value = -28193

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 89286840
$ws.Range("I4").Value = 156250620
$ws.Range("K4").Value = 468751860
$ws.Range("M4").Value = -468751748
$ws.Range("H23").Value = 67.53333000000001
$ws.Range("J23").Value = 136
$ws.Range("L23").Value = 408
$ws.Range("N23").Value = -878
$ws.Range("H34").Value = 2179.7273
$ws.Range("I34").Value = 867.6
$ws.Range("J34").Value = 3273.1667
$ws.Range("K34").Value = 2602.8
$ws.Range("L34").Value = 9819.500100000001
$ws.Range("M34").Value = -2518.8
$ws.Range("N34").Value = -9987.500100000001
$ws.Range("H39").Value = 3749
$ws.Range("J39").Value = 3999.8
$ws.Range("L39").Value = 11999.4
$ws.Range("N39").Value = -12587.4
$ws.Range("H51").Value = 198.66667
$ws.Range("J51").Value = 198.66667
$ws.Range("L51").Value = 596.00001
$ws.Range("N51").Value = -1516.00001
$ws.Range("H118").Value = 3169.4614
$ws.Range("I118").Value = 3567
$ws.Range("J118").Value = 2275
$ws.Range("K118").Value = 10701
$ws.Range("L118").Value = 6825
$ws.Range("M118").Value = -9458
$ws.Range("N118").Value = -9311
$ws.Range("H131").Value = 700705.6
$ws.Range("I131").Value = 795.82355
$ws.Range("J131").Value = 2022757.5
$ws.Range("K131").Value = 2387.47065
$ws.Range("L131").Value = 6068272.5
$ws.Range("M131").Value = 2652.52935
$ws.Range("N131").Value = -6078352.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 548
$ws.Range("I107").Value = 548
$ws.Range("K107").Value = 548
$ws.Range("M107").Value = 1372
$ws.Range("H126").Value = 5998
$ws.Range("J126").Value = 7998.25
$ws.Range("L126").Value = 23994.75
$ws.Range("N126").Value = -28934.75
$ws.Range("H132").Value = 7438.625
$ws.Range("I132").Value = 6881.769
$ws.Range("J132").Value = 8096.727
$ws.Range("K132").Value = 20645.307
$ws.Range("L132").Value = 24290.181
$ws.Range("M132").Value = -18115.307
$ws.Range("N132").Value = -29350.181
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -95100

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 891.06665
$ws.Range("I22").Value = 461.61905
$ws.Range("J22").Value = 1893.1111
$ws.Range("K22").Value = 461.61905
$ws.Range("L22").Value = 1893.1111
$ws.Range("M22").Value = -166.61905
$ws.Range("N22").Value = -2483.1111
$ws.Range("H27").Value = 891.06665
$ws.Range("I27").Value = 461.61905
$ws.Range("J27").Value = 1893.1111
$ws.Range("K27").Value = 461.61905
$ws.Range("L27").Value = 1893.1111
$ws.Range("M27").Value = -354.61905
$ws.Range("N27").Value = -2107.1111
$ws.Range("H40").Value = 6034.8184
$ws.Range("I40").Value = 4709.222
$ws.Range("K40").Value = 4709.222
$ws.Range("M40").Value = -4573.222
$ws.Range("H55").Value = 968.86365
$ws.Range("I55").Value = 604.25
$ws.Range("J55").Value = 1177.2142
$ws.Range("K55").Value = 604.25
$ws.Range("L55").Value = 1177.2142
$ws.Range("M55").Value = -431.25
$ws.Range("N55").Value = -1523.2142
$ws.Range("H68").Value = 6964.4136
$ws.Range("J68").Value = 6242.4287
$ws.Range("L68").Value = 6242.4287
$ws.Range("N68").Value = -7740.4287
$ws.Range("H71").Value = 6964.4136
$ws.Range("J71").Value = 6242.4287
$ws.Range("L71").Value = 31212.1435
$ws.Range("N71").Value = -38700.14350000001
$ws.Range("H122").Value = 4600.9443
$ws.Range("J122").Value = 4699.5
$ws.Range("L122").Value = 14098.5
$ws.Range("N122").Value = -18998.5
$ws.Range("H132").Value = 4441.811
$ws.Range("I132").Value = 3736.2666
$ws.Range("K132").Value = 11208.7998
$ws.Range("M132").Value = -8678.799800000001
$ws.Range("H136").Value = 5086.4443
$ws.Range("I136").Value = 4357.5454
$ws.Range("K136").Value = 13072.6362
$ws.Range("M136").Value = -10522.6362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 45000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 45000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 45000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -45586
$ws.Range("H28").Value = 500
$ws.Range("J28").Value = 500
$ws.Range("L28").Value = 500
$ws.Range("N28").Value = -1196
$ws.Range("H96").Value = 22534.273
$ws.Range("J96").Value = 26708.666
$ws.Range("L96").Value = 26708.666
$ws.Range("N96").Value = -29454.666
$ws.Range("H100").Value = 486.57144
$ws.Range("I100").Value = 492.66666
$ws.Range("K100").Value = 985.33332
$ws.Range("M100").Value = -444.33332
$ws.Range("H122").Value = 3384.0417
$ws.Range("I122").Value = 2458.6875
$ws.Range("J122").Value = 5234.75
$ws.Range("K122").Value = 7376.0625
$ws.Range("L122").Value = 15704.25
$ws.Range("M122").Value = -4926.0625
$ws.Range("N122").Value = -20604.25
